$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

$ws.Cells.Item($row, 1).Value = 38
$ws.Cells.Item($row, 2).Value = "armenia"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45189.5
$ws.Cells.Item($row, 6).Value = "Alashkert"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Ararat Yerevan"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.48
$ws.Cells.Item($row, 11).Value = "18/09/2023 23:12"
$ws.Cells.Item($row, 12).Value = 1.32
$ws.Cells.Item($row, 13).Value = "20/09/2023 11:53"
$ws.Cells.Item($row, 14).Value = 4.02
$ws.Cells.Item($row, 15).Value = "18/09/2023 23:12"
$ws.Cells.Item($row, 16).Value = 5.11
$ws.Cells.Item($row, 17).Value = "20/09/2023 11:59"
$ws.Cells.Item($row, 18).Value = 5.97
$ws.Cells.Item($row, 19).Value = "18/09/2023 23:12"
$ws.Cells.Item($row, 20).Value = 9.880000000000001
$ws.Cells.Item($row, 21).Value = "20/09/2023 11:59"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/alashkert-ararat-yerevan/2sxoU7FM/"

# Copy styles from row 38 to new row 39 (cell A gets s="1", cell E gets s="2")
$ws.Cells.Item(38, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(38, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)  # xlPasteFormats
